$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Recorded window resize (bookViews/workbookView windowHeight 12330 -> 13770,
# i.e. 616.5pt -> 688.5pt) from the original editing session.
$excel.ActiveWindow.Height = 688.5

$ws.Range("E2").Value = 312.25115369684761
$ws.Range("F2").Value = -8.0318491234546219
$ws.Range("G2").Value = [double]"-9.6815087654537724E-2"
$ws.Range("H2").Value = -10.112054942940608
$ws.Range("I2").Value = [double]"-1.226178284972379E-2"
$ws.Range("J2").Value = [double]"5.0213410604839082E-2"
$ws.Range("K2").Value = -0.42205169102136431
$ws.Range("L2").Value = -0.64448707724465892
$ws.Range("M2").Value = -0.55397777535873804
$ws.Range("N2").Value = [double]"9.514535503055277E-2"
$ws.Range("O2").Value = -0.34631943650294505
$ws.Range("P2").Value = 15
$ws.Range("Q2").Value = 1
$ws.Range("E3").Value = 331.28137564737455
$ws.Range("F3").Value = -8.1305200810389877
$ws.Range("G3").Value = [double]"-8.694799189610114E-2"
$ws.Range("H3").Value = -10.226541325427219
$ws.Range("I3").Value = [double]"-2.0852155175806741E-3"
$ws.Range("J3").Value = [double]"1.1014055162195335E-2"
$ws.Range("K3").Value = -0.49287020697214112
$ws.Range("L3").Value = -0.62678244825696472
$ws.Range("M3").Value = -0.57960059818927179
$ws.Range("N3").Value = [double]"2.1698777989694884E-2"
$ws.Range("O3").Value = -0.39674929451277063
$ws.Range("E4").Value = 336.24534636048384
$ws.Range("F4").Value = -8.1316532899987468
$ws.Range("G4").Value = [double]"-8.6834671000125319E-2"
$ws.Range("H4").Value = -10.225161973110446
$ws.Range("I4").Value = [double]"-2.2078246124047629E-3"
$ws.Range("J4").Value = [double]"-3.9082174032528672E-3"
$ws.Range("K4").Value = -0.51523711510389569
$ws.Range("L4").Value = -0.62119072122402608
$ws.Range("M4").Value = -0.58580749662443465
$ws.Range("N4").Value = [double]"6.7350616961001641E-3"
$ws.Range("O4").Value = -0.40565821981078953
$ws.Range("E5").Value = 338.06386436367933
$ws.Range("F5").Value = -8.1320724055571301
$ws.Range("G5").Value = [double]"-8.6792759444286882E-2"
$ws.Range("H5").Value = -10.224661707905506
$ws.Range("I5").Value = [double]"-2.2522926306216551E-3"
$ws.Range("J5").Value = [double]"-9.2643405795790335E-3"
$ws.Range("K5").Value = -0.52343039232726296
$ws.Range("L5").Value = -0.61914240191818426
$ws.Range("M5").Value = -0.58803573528058151
$ws.Range("N5").Value = [double]"1.362725487747074E-3"
$ws.Range("O5").Value = -0.40885699416256294
$ws.Range("E6").Value = 338.48756076366965
$ws.Range("F6").Value = -8.1321697893086604
$ws.Range("G6").Value = [double]"-8.6783021069133937E-2"
$ws.Range("H6").Value = -10.224544813361623
$ws.Range("I6").Value = [double]"-2.2626832567445021E-3"
$ws.Range("J6").Value = [double]"-1.0504055472598472E-2"
$ws.Range("K6").Value = -0.52533938619361376
$ws.Range("L6").Value = -0.61866515345159656
$ws.Range("M6").Value = -0.58855145362523897
$ws.Range("N6").Value = [double]"1.1934618735454094E-4"
$ws.Range("O6").Value = -0.40959730489014479
$ws.Range("E7").Value = 338.5290405049563
$ws.Range("F7").Value = -8.1321792221476752
$ws.Range("G7").Value = [double]"-8.6782077785232503E-2"
$ws.Range("H7").Value = -10.224533241035418
$ws.Range("I7").Value = [double]"-2.2637119079628176E-3"
$ws.Range("J7").Value = [double]"-1.0625278644965541E-2"
$ws.Range("K7").Value = -0.52552629319340216
$ws.Range("L7").Value = -0.61861842670164946
$ws.Range("M7").Value = -0.58860187291248312
$ws.Range("N7").Value = [double]"-2.20156157948459E-6"
$ws.Range("O7").Value = -0.40966966821166217
$ws.Range("E8").Value = 338.52440651440168
$ws.Range("F8").Value = -8.1321781693115582
$ws.Range("G8").Value = [double]"-8.678218306884411E-2"
$ws.Range("H8").Value = -10.224534535096094
$ws.Range("I8").Value = [double]"-2.2635968803472739E-3"
$ws.Range("J8").Value = [double]"-1.0611737220627471E-2"
$ws.Range("K8").Value = -0.52550541234683923
$ws.Range("L8").Value = -0.61862364691329019
$ws.Range("M8").Value = -0.58859624083607764
$ws.Range("N8").Value = [double]"1.1375795425205837E-5"
$ws.Range("O8").Value = -0.40966158500906424
$ws.Range("E9").Value = 338.53162810627509
$ws.Range("F9").Value = -8.1321796565655156
$ws.Range("G9").Value = [double]"-8.6782034343448444E-2"
$ws.Range("H9").Value = -10.224532323325988
$ws.Range("I9").Value = [double]"-2.2637934821344574E-3"
$ws.Range("J9").Value = [double]"-1.0632873780386021E-2"
$ws.Range("K9").Value = -0.52553797904725341
$ws.Range("L9").Value = -0.61861550523818665
$ws.Range("M9").Value = -0.58860501777509322
$ws.Range("N9").Value = [double]"-9.7655676213337372E-6"
$ws.Range("O9").Value = -0.40967416156549408
$ws.Range("E10").Value = 338.5283593748498
$ws.Range("F10").Value = -8.1321789139142879
$ws.Range("G10").Value = [double]"-8.6782108608571162E-2"
$ws.Range("H10").Value = -10.224533236132462
$ws.Range("I10").Value = [double]"-2.263712343781088E-3"
$ws.Range("J10").Value = [double]"-1.062332209442074E-2"
$ws.Range("K10").Value = -0.52552325008318279
$ws.Range("L10").Value = -0.6186191874792043
$ws.Range("M10").Value = -0.58860104508907696
$ws.Range("N10").Value = [double]"-1.8853538563678995E-7"
$ws.Range("O10").Value = -0.40966845993349155
$ws.Range("E11").Value = 338.5283593748498
$ws.Range("F11").Value = -8.1321789139142879
$ws.Range("G11").Value = [double]"-8.6782108608571162E-2"
$ws.Range("H11").Value = -10.224533236132462
$ws.Range("I11").Value = [double]"-2.263712343781088E-3"
$ws.Range("J11").Value = [double]"-1.062332209442074E-2"
$ws.Range("K11").Value = -0.52552325008318279
$ws.Range("L11").Value = -0.6186191874792043
$ws.Range("M11").Value = -0.58860104508907696
$ws.Range("N11").Value = [double]"-1.8853538563678995E-7"
$ws.Range("O11").Value = -0.40966845993349155
